{"js": "// \"Better handle copying paragraph styles\": inside each text table cell,\n// a join-marker paragraph styled MSC_Join is immediately followed by an\n// \"[...]\" placeholder paragraph and a blank paragraph that belong to the\n// same join group, but those two paragraphs were left without an explicit\n// style. Copy the MSC_Join paragraph style onto those two sibling\n// paragraphs so the whole group is consistently styled.\n//\n// Work cell-by-cell (each table cell's own paragraph collection), so \"the\n// next two paragraphs\" never crosses a table-cell boundary. (Note:\n// body.paragraphs is a flattened view that also includes paragraphs that\n// live inside table cells, so it is not used here for locating siblings.)\n\nconst JOIN_STYLE = \"MSC_Join\";\n\nfunction applyJoinStyleForward(items) {\n  // Snapshot the styles that existed *before* any mutation so that a\n  // paragraph we just switched to MSC_Join in this pass is never treated\n  // as if it were already an (original) join-marker paragraph -- that\n  // would incorrectly cascade the style onto further, unrelated\n  // paragraphs.\n  const originalStyles = items.map((p) => p.style);\n\n  for (let i = 0; i < items.length; i++) {\n    if (originalStyles[i] === JOIN_STYLE) {\n      for (let offset = 1; offset <= 2 && i + offset < items.length; offset++) {\n        const target = items[i + offset];\n        if (target.style !== JOIN_STYLE) {\n          target.style = JOIN_STYLE;\n        }\n      }\n    }\n  }\n}\n\nconst body = context.document.body;\n\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Collect every table cell's own paragraph collection.\nconst cellParagraphCollections = [];\nfor (let t = 0; t < tables.items.length; t++) {\n  const rows = tables.items[t].rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (let r = 0; r < rows.items.length; r++) {\n    const cells = rows.items[r].cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    for (let c = 0; c < cells.items.length; c++) {\n      const cellParagraphs = cells.items[c].body.paragraphs;\n      cellParagraphs.load(\"items/style\");\n      cellParagraphCollections.push(cellParagraphs);\n    }\n  }\n}\n\nawait context.sync();\n\n// Apply within each table cell's own paragraph list (container-scoped).\nfor (const cellParagraphs of cellParagraphCollections) {\n  applyJoinStyleForward(cellParagraphs.items);\n}\n\nawait context.sync();\n", "ps1": "# \"Better handle copying paragraph styles\": a join-marker paragraph styled\n# MSC_Join is immediately followed (inside the same table cell) by an\n# \"[...]\" placeholder paragraph and a blank paragraph that belong to the\n# same join group, but those two paragraphs were left without an explicit\n# style. Copy the MSC_Join paragraph style onto those two sibling\n# paragraphs so the whole group is consistently styled.\n\n$d = $word.ActiveDocument\n\n$JoinStyleName = \"MSC_Join\"\n\n$paragraphs = $d.Paragraphs\n$count = $paragraphs.Count\n\n# Snapshot every paragraph's style *before* making any change. Using a\n# snapshot (rather than re-reading live style after each write) keeps a\n# paragraph we just switched to MSC_Join in this pass from being treated\n# as if it were already an original join-marker paragraph -- that would\n# incorrectly cascade the style onto further, unrelated paragraphs.\n$originalStyles = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $originalStyles += $paragraphs.Item($i).Range.Style.NameLocal\n}\n\nfor ($i = 1; $i -le $count; $i++) {\n    if ($originalStyles[$i - 1] -eq $JoinStyleName) {\n        $maxOffset = 2\n        if ($count - $i -lt $maxOffset) {\n            $maxOffset = $count - $i\n        }\n        for ($offset = 1; $offset -le $maxOffset; $offset++) {\n            $target = $paragraphs.Item($i + $offset)\n            if ($target.Range.Style.NameLocal -ne $JoinStyleName) {\n                $target.Range.Style = $JoinStyleName\n            }\n        }\n    }\n}\n"}
